# Update "想去人数" (number of people wanting to go) counts for a handful
# of events across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1526
$ws1.Range("F5").Value = 789
$ws1.Range("F7").Value = 13259
$ws1.Range("F13").Value = 690
$ws1.Range("F20").Value = 277

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 63
$ws2.Range("F9").Value = 34

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1526
$ws4.Range("F7").Value = 789
$ws4.Range("F9").Value = 13259
$ws4.Range("F15").Value = 690
$ws4.Range("F24").Value = 63
$ws4.Range("F27").Value = 277
$ws4.Range("F34").Value = 34
